$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting existing rows 8-20 down to 9-21.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44797
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11500
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 460
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
